$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the English header from "en-gb" to "en_gb" ---
$ws.Range("B1").Value = "en_gb"

# --- Add the three new "creation.*" translation rows (7-10) ---

# Row 7: creation.placeholder
$ws.Range("A7").Value = "creation.placeholder"
$ws.Range("B7").Value = "Mario playing with…"
$ws.Range("C7").Value = "Mario jouant avec..."
$ws.Range("D7").Value = "Mario spielt mit..."
$ws.Range("E7").Value = "Mario gioca con..."
$ws.Range("F7").Value = "Mario jugando con..."

# Row 8: creation.option1
$ws.Range("A8").Value = "creation.option1"
$ws.Range("B8").Value = "Better quality (increased price)"
$ws.Range("C8").Value = "Meilleure qualité d'image (augmentation du prix)"
$ws.Range("D8").Value = "Bessere Bildqualität (Preiserhöhung)"
$ws.Range("E8").Value = "Migliore qualità dell'immagine (aumento del prezzo)"
$ws.Range("F8").Value = "Mejor calidad de imagen (aumento de precio)"

# Row 9: creation.option2
$ws.Range("A9").Value = "creation.option2"
$ws.Range("B9").Value = "More difficult colouring"
$ws.Range("C9").Value = "Coloriage plus difficile"
$ws.Range("D9").Value = "Schwierigeres Ausmalen"
$ws.Range("E9").Value = "Colorazione più difficile"
$ws.Range("F9").Value = "Colorear más difícil"

# Row 10: creation.option3
$ws.Range("A10").Value = "creation.option3"
$ws.Range("B10").Value = "Only one page"
$ws.Range("C10").Value = "Une seule page"
$ws.Range("D10").Value = "Nur eine Seite"
$ws.Range("E10").Value = "Solo una pagina"
$ws.Range("F10").Value = "Solo una página"

# --- Update the view: scroll back to top-left A1 and select B10 ---
$aw = $excel.ActiveWindow
$aw.ScrollColumn = 1
$aw.ScrollRow = 1
$ws.Range("B10").Select()
